$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the weekly records (rows 2-4) so that:
#   row2 <- old row3, row3 <- old row4, row4 <- old row2
# Columns D, J, K, L, M, O, P hold the values that change; capture the
# originals first so the rotation can be applied without clobbering data.
$D2 = $ws.Range("D2").Value2
$D3 = $ws.Range("D3").Value2
$D4 = $ws.Range("D4").Value2

$J2 = $ws.Range("J2").Value2
$J3 = $ws.Range("J3").Value2
$J4 = $ws.Range("J4").Value2

$K2 = $ws.Range("K2").Value2
$K3 = $ws.Range("K3").Value2
$K4 = $ws.Range("K4").Value2

$L2 = $ws.Range("L2").Value2
$L3 = $ws.Range("L3").Value2
$L4 = $ws.Range("L4").Value2

$M2 = $ws.Range("M2").Value2
$M3 = $ws.Range("M3").Value2
$M4 = $ws.Range("M4").Value2

$O2 = $ws.Range("O2").Value2
$O3 = $ws.Range("O3").Value2
$O4 = $ws.Range("O4").Value2

$P2 = $ws.Range("P2").Value2
$P3 = $ws.Range("P3").Value2
$P4 = $ws.Range("P4").Value2

$ws.Range("D2").Value = $D3
$ws.Range("D3").Value = $D4
$ws.Range("D4").Value = $D2

$ws.Range("J2").Value = $J3
$ws.Range("J3").Value = $J4
$ws.Range("J4").Value = $J2

$ws.Range("K2").Value = $K3
$ws.Range("K3").Value = $K4
$ws.Range("K4").Value = $K2

$ws.Range("L2").Value = $L3
$ws.Range("L3").Value = $L4
$ws.Range("L4").Value = $L2

$ws.Range("M2").Value = $M3
$ws.Range("M3").Value = $M4
$ws.Range("M4").Value = $M2

$ws.Range("O2").Value = $O3
$ws.Range("O3").Value = $O4
$ws.Range("O4").Value = $O2

$ws.Range("P2").Value = $P3
$ws.Range("P3").Value = $P4
$ws.Range("P4").Value = $P2

$wb.Save()
